# Auto-generated Excel COM-interop edit script
# Applies the "Updated cryptos list" data refresh (coin prices / 1h volume %) to Sheet1,
# matching the GitHub Actions commit that refreshed cryptos.xlsx.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '60.372.72'
$ws.Range('E2').Value = '  -3.58%  '
# Row 3
$ws.Range('D3').Value = '3.308.40'
$ws.Range('E3').Value = '  -4.03%  '
# Row 4
$ws.Range('E4').Value = '  -0.01%  '
# Row 5
$ws.Range('D5').Value = "'560.10"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.10%  '
# Row 6
$ws.Range('D6').Value = "'144.12"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.04%  '
# Row 7
$ws.Range('E7').Value = '  -0.04%  '
# Row 8
$ws.Range('D8').Value = '3.311.49'
$ws.Range('E8').Value = '  -3.88%  '
# Row 9
$ws.Range('E9').Value = '  -0.50%  '
# Row 10
$ws.Range('D10').Value = "'7.82"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.72%  '
# Row 11
$ws.Range('E11').Value = '  -3.31%  '
# Row 12
$ws.Range('D12').Value = "'0.408"
$ws.Range('D12').Style = 'Normal'
# Row 13
$ws.Range('D13').Value = '3.873.26'
$ws.Range('E13').Value = '  -4.11%  '
# Row 15
$ws.Range('D15').Value = "'27.26"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.72%  '
# Row 16
$ws.Range('D16').Value = '3.317.92'
$ws.Range('E16').Value = '  -3.76%  '
# Row 17
$ws.Range('E17').Value = '  -3.47%  '
# Row 18
$ws.Range('D18').Value = '60.367.51'
$ws.Range('E18').Value = '  -3.71%  '
# Row 19
$ws.Range('D19').Value = "'6.15"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.00%  '
# Row 20
$ws.Range('D20').Value = "'14.39"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.46%  '
# Row 21
$ws.Range('D21').Value = "'8.59"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.44%  '
# Row 22
$ws.Range('D22').Value = "'373.46"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.45%  '
# Row 23
$ws.Range('D23').Value = "'74.14"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.46%  '
# Row 24
$ws.Range('E24').Value = '  -3.12%  '
# Row 25
$ws.Range('E25').Value = '  -0.06%  '
# Row 26
$ws.Range('D26').Value = '3.469.90'
$ws.Range('E26').Value = '  -3.21%  '
# Row 28
$ws.Range('D28').Value = "'0.173"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.56%  '
# Row 29
$ws.Range('D29').Value = "'0.997"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.21%  '
# Row 30
$ws.Range('D30').Value = "'7.22"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.83%  '
# Row 31
$ws.Range('E31').Value = '  +0.04%  '
# Row 32
$ws.Range('E32').Value = '  -3.18%  '
# Row 33
$ws.Range('D33').Value = "'7.62"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.73%  '
# Row 34
$ws.Range('D34').Value = "'22.61"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.64%  '
# Row 35
$ws.Range('E35').Value = '  -4.60%  '
# Row 36
$ws.Range('D36').Value = "'5.18"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.56%  '
# Row 37
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').Value = "'6.76"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.68%  '
# Row 38
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').Value = "'165.90"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.94%  '
# Row 39
$ws.Range('E39').Value = '  -6.32%  '
# Row 40
$ws.Range('D40').Value = "'27.75"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -13.38%  '
# Row 41
$ws.Range('D41').Value = '3.339.70'
$ws.Range('E41').Value = '  -4.11%  '
# Row 42
$ws.Range('D42').Value = "'0.0739"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.74%  '
# Row 43
$ws.Range('D43').Value = "'41.96"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.76%  '
# Row 44
$ws.Range('D44').Value = "'0.753"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.16%  '
# Row 45
$ws.Range('E45').Value = '  -3.96%  '
# Row 46
$ws.Range('E46').Value = '  -4.78%  '
# Row 47
$ws.Range('E47').Value = '  -4.19%  '
# Row 48
$ws.Range('D48').Value = '2.376.12'
$ws.Range('E48').Value = '  -7.45%  '
# Row 49
$ws.Range('E49').Value = '  -0.09%  '
# Row 50
$ws.Range('E50').Value = '  -4.50%  '
# Row 51
$ws.Range('D51').Value = "'21.68"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.76%  '
